$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.27%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.21%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.946"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.79%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07144"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-9.19%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-12.78%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.628"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.73%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.736"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8966"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.22%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1647"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.65%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07629"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.25%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08051"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.08%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03047"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.71%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1002"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.41%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001487"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.68%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005768"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.45%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.467"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.07%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.71%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3276"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.00%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1312"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.38%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.047"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.38%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.68%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04510"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.59%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.13%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.22%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001250"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.05%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01610"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.11%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04361"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-9.29%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007387"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.40%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1303"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.46%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002029"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-14.01%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009160"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.18%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005967"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.67%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.00%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.246"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "172.73%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.23%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
